# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" after "2021-Q4" (and before "总计"),
#    populated with the per-fund holdings for that quarter.
# 2. Insert a new leading row into the "总计" (totals) summary sheet for
#    2022-Q1 (11 holdings, 0.38 亿元 total value), shifting the existing
#    rows down and renumbering the index column.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# --- 1. Create the "2022-Q1" sheet, positioned after "2021-Q4" -------------
$afterSheet = $sheets.Item("2021-Q4")
$newSheet = $sheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "2022-Q1"

# Reuse the header formatting (bold/centered/bordered) from an existing
# per-fund sheet so the new sheet matches the workbook's look.
$template = $sheets.Item("2021-Q4")
$template.Range("B1:H1").Copy($newSheet.Range("B1:H1"))
$template.Range("A2").Copy($newSheet.Range("A2:A12"))

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Fund code / name / scale / position columns are stored as text (codes like
# "012010" must keep their leading zero), so force text format before the
# bulk write.
$newSheet.Range("B2:G12").NumberFormat = "@"

$fundData = @"
0	510081	长盛动态精选混合	3.15	60.76	3.61	0.1137	5
1	012010	富国泰享回报6个月持有期混合型证券投资基金A	9.29	29.91	1.16	0.1078	6
2	004945	长信中证500指数增强	2.81	92.77	1.82	0.0511	4
3	005251	银华多元动力灵活配置混合	1.74	93.75	2.32	0.0404	8
4	009726	招商中证500等权重指数增强A	1.87	91.11	1.48	0.0277	5
5	005357	富国国企改革灵活配置混合	1.13	87.21	2.14	0.0242	8
6	009727	招商中证500等权重指数增强C	0.69	91.11	1.48	0.0102	5
7	003670	中融物联网主题灵活配置混合	0.15	79.61	3.21	0.0048	9
8	515510	嘉实中证500成长估值ETF	0.15	98.79	1.24	0.0019	6
9	001744	诺安进取回报灵活配置混合	0.04	62.10	4.62	0.0018	4
10	012011	富国泰享回报6个月持有期混合型证券投资基金C	0.09	29.91	1.16	0.0010	6
"@

$fundLines = $fundData -split "`n"
$fundRowCount = $fundLines.Count
$fundArr = New-Object 'object[,]' $fundRowCount,8
for ($i = 0; $i -lt $fundRowCount; $i++) {
    $cols = $fundLines[$i] -split "`t"
    for ($j = 0; $j -lt 8; $j++) {
        $fundArr[$i, $j] = $cols[$j]
    }
}
$newSheet.Range("A2:H12").Value = $fundArr

# --- 2. Add the 2022-Q1 summary row to "总计" -------------------------------
$totalSheet = $sheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The blank row Insert() leaves behind borrows the header's bold style for
# B:D; reset those to plain/default, then restore the index column's
# (column A) styling from the row it pushed down.
$totalSheet.Range("B2:D2").Style = "Normal"
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 11
$totalSheet.Range("D2").Value = 0.38

# Renumber the index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3

# Restore the originally active sheet/selection (adding sheets moves focus).
$sheets.Item("2021-Q1").Activate()
[void]$sheets.Item("2021-Q1").Range("A1").Select()
